$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the next diagonal of "_diff" error values (ifo GDP component analysis preprocessing)
$ws.Range("K15").Value = 0.1683237681281231
$ws.Range("J16").Value = 0.1722916656412322
$ws.Range("I17").Value = 0.3079317558114735
$ws.Range("H18").Value = 0.06712557395580883
$ws.Range("G19").Value = 0.02179435870371246
$ws.Range("F20").Value = -0.04506706323234141
$ws.Range("E21").Value = -0.07465326558905801
$ws.Range("D22").Value = -0.0928039223186989
$ws.Range("C23").Value = -0.1108357465673982
$ws.Range("B24").Value = -0.1624199859130616
